$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last existing data row (149) down across the new rows so that the new rows
# inherit the same cell styles/number formats (bold+border on col A, date format on col E)
$ws.Range("A149:AC149").Copy($ws.Range("A150:AC154"))

# ---- Row 150 ----
$ws.Range("A150").Value = 148
$ws.Range("B150").Value = 6957486
$ws.Range("C150").Value = "Slovenia Prva Liga"
$ws.Range("D150").Value = "Slovenia Prva Liga"
$ws.Range("E150").Value = 45396.63541666666
$ws.Range("F150").Value = "NK Domzale"
$ws.Range("G150").Value = "Olimpija Ljubljana"
$ws.Range("H150").Value = 1
$ws.Range("I150").Value = 3
$ws.Range("J150").Value = "A"
$ws.Range("K150").Value = 4.75
$ws.Range("L150").Value = 4
$ws.Range("M150").Value = 1.571
$ws.Range("N150").Value = 4.333
$ws.Range("O150").Value = 3.8
$ws.Range("P150").Value = 1.65
$ws.Range("Q150").Value = 0.75
$ws.Range("R150").Value = 1.95
$ws.Range("S150").Value = 1.85
$ws.Range("T150").Value = 2.75
$ws.Range("U150").Value = 1.9
$ws.Range("V150").Value = 1.9
$ws.Range("W150").Value = -1
$ws.Range("X150").Value = -1
$ws.Range("Y150").Value = 0.6499999999999999
$ws.Range("Z150").Value = -1
$ws.Range("AA150").Value = 0.8500000000000001
$ws.Range("AB150").Value = 0.8999999999999999
$ws.Range("AC150").Value = -1

# ---- Row 151 ----
$ws.Range("A151").Value = 149
$ws.Range("B151").Value = 6961772
$ws.Range("C151").Value = "Slovenia Prva Liga"
$ws.Range("D151").Value = "Slovenia Prva Liga"
$ws.Range("E151").Value = 45397.52083333334
$ws.Range("F151").Value = "NK Aluminij"
$ws.Range("G151").Value = "NK Radomlje"
$ws.Range("H151").Value = 1
$ws.Range("I151").Value = 1
$ws.Range("J151").Value = "D"
$ws.Range("K151").Value = 3.2
$ws.Range("L151").Value = 3.25
$ws.Range("M151").Value = 2.1
$ws.Range("N151").Value = 3.4
$ws.Range("O151").Value = 3.3
$ws.Range("P151").Value = 2
$ws.Range("Q151").Value = 0.25
$ws.Range("R151").Value = 1.975
$ws.Range("S151").Value = 1.825
$ws.Range("T151").Value = 2.5
$ws.Range("U151").Value = 1.925
$ws.Range("V151").Value = 1.875
$ws.Range("W151").Value = -1
$ws.Range("X151").Value = 2.3
$ws.Range("Y151").Value = -1
$ws.Range("Z151").Value = 0.4875
$ws.Range("AA151").Value = -0.5
$ws.Range("AB151").Value = -1
$ws.Range("AC151").Value = 0.875

# ---- Row 152 ----
$ws.Range("A152").Value = 150
$ws.Range("B152").Value = 8100619
$ws.Range("C152").Value = "Slovenia Prva Liga"
$ws.Range("D152").Value = "Slovenia Prva Liga"
$ws.Range("E152").Value = 45400.45833333334
$ws.Range("F152").Value = "NK Domzale"
$ws.Range("G152").Value = "FC Koper"
$ws.Range("K152").Value = 3.1
$ws.Range("L152").Value = 3.5
$ws.Range("M152").Value = 2.05
$ws.Range("N152").Value = 2.9
$ws.Range("O152").Value = 3.5
$ws.Range("P152").Value = 2.2
$ws.Range("Q152").Value = 0.25
$ws.Range("R152").Value = 1.85
$ws.Range("S152").Value = 1.95
$ws.Range("T152").Value = 2.75
$ws.Range("U152").Value = 1.925
$ws.Range("V152").Value = 1.875
$ws.Range("W152").Value = 0
$ws.Range("X152").Value = 0
$ws.Range("Y152").Value = 0
$ws.Range("Z152").Value = 0
$ws.Range("AA152").Value = 0
$ws.Range("H152:J152").ClearContents()
$ws.Range("AB152:AC152").ClearContents()

# ---- Row 153 ----
$ws.Range("A153").Value = 151
$ws.Range("B153").Value = 6994887
$ws.Range("C153").Value = "Slovenia Prva Liga"
$ws.Range("D153").Value = "Slovenia Prva Liga"
$ws.Range("E153").Value = 45402.41666666666
$ws.Range("F153").Value = "NK Radomlje"
$ws.Range("G153").Value = "NS Mura"
$ws.Range("K153").Value = 2.625
$ws.Range("L153").Value = 3.25
$ws.Range("M153").Value = 2.45
$ws.Range("N153").Value = 2.25
$ws.Range("O153").Value = 3.25
$ws.Range("P153").Value = 2.8
$ws.Range("Q153").Value = -0.25
$ws.Range("R153").Value = 2.025
$ws.Range("S153").Value = 1.775
$ws.Range("T153").Value = 2.25
$ws.Range("U153").Value = 1.8
$ws.Range("V153").Value = 2
$ws.Range("W153").Value = 0
$ws.Range("X153").Value = 0
$ws.Range("Y153").Value = 0
$ws.Range("Z153").Value = 0
$ws.Range("AA153").Value = 0
$ws.Range("H153:J153").ClearContents()
$ws.Range("AB153:AC153").ClearContents()

# ---- Row 154 ----
$ws.Range("A154").Value = 152
$ws.Range("B154").Value = 6998172
$ws.Range("C154").Value = "Slovenia Prva Liga"
$ws.Range("D154").Value = "Slovenia Prva Liga"
$ws.Range("E154").Value = 45402.52083333334
$ws.Range("F154").Value = "Olimpija Ljubljana"
$ws.Range("G154").Value = "NK Aluminij"
$ws.Range("K154").Value = 1.25
$ws.Range("L154").Value = 5.25
$ws.Range("M154").Value = 9.5
$ws.Range("N154").Value = 1.3
$ws.Range("O154").Value = 5
$ws.Range("P154").Value = 8
$ws.Range("Q154").Value = -1.5
$ws.Range("R154").Value = 1.85
$ws.Range("S154").Value = 1.95
$ws.Range("T154").Value = 3
$ws.Range("U154").Value = 1.9
$ws.Range("V154").Value = 1.9
$ws.Range("W154").Value = 0
$ws.Range("X154").Value = 0
$ws.Range("Y154").Value = 0
$ws.Range("Z154").Value = 0
$ws.Range("AA154").Value = 0
$ws.Range("H154:J154").ClearContents()
$ws.Range("AB154:AC154").ClearContents()

